$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Nathan Weiler"
$ws.Range("B4").Value = "nweiler@uoguelph.ca"
